# Apply Karthik's output edits to the REPORT worksheet.
# Summary of the change:
#  - Row 29: remove the stray empty (but date-styled) cell D29.
#  - Row 39: this row was an accidental duplicate of the "Selection 3" block
#    header (B:G) with a wrong count; clear B39:G39 so only the H39 note remains.
#  - Rows 44-94: delete the duplicated trailing blocks that were left over
#    from copy/pasting the "Selection" sections multiple times.
#  - Column D: the width was manually resized (no longer auto "best fit").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REPORT")

# Remove the leftover empty styled cell in D29.
$ws.Range("D29").Clear()

# Row 39 had an erroneous duplicate block header (B:G); strip it back down
# to just the trailing note in column H.
$ws.Range("B39:G39").Clear()

# Delete the duplicated tail blocks (rows 44-94) entirely, shifting
# everything below up and shrinking the sheet's used range to A26:H43.
$ws.Rows("44:94").Delete()

# Column D was manually resized (bestFit flag cleared) to a slightly wider,
# fixed width.
$ws.Columns.Item(4).ColumnWidth = 9.5
